$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------------
# Data for the "Journal de dev" tracking sheet (A1:D22)
# Columns: A=Classe, B=Fonction, C=Développeur, D=Etat
#
# Cell values are keyed by address so we can fill the sheet in the same
# "block by block" order the sheet was originally authored in (each
# merged-row group was filled column by column: A, then B, then C, then
# D) which reproduces the exact shared-string table ordering.
# ----------------------------------------------------------------------
$values = @{
    "A1" = "Classe";  "B1" = "Fonction"; "C1" = "Développeur"; "D1" = "Etat"
    "A2" = "Argent";                                           "D2" = "Fait"
    "A3" = "Boss";                                              "D3" = "Fait"
    "A4" = "Gangster";                                          "D4" = "Fait"
    "A5" = "Case";                                              "D5" = "Fait"
                      "B6" = "Affichage";                       "D6" = "Fait"
    "A7" = "Dé";                                                "D7" = "Fait"
                      "B8" = "Résultat";                        "D8" = "Fait"
    "A9" = "Plateau";                                           "D9" = "A faire"
                      "B10" = "Afficher";                       "D10" = "Fait"
                      "B11" = "Remplissage";                    "D11" = "A faire"
    "A12" = "Prison";                                           "D12" = "Fait"
                      "B13" = "Affichage";                      "D13" = "Fait"
    "A14" = "Joueur";                                           "D14" = "Fait"
                      "B15" = "Déplacement";                    "D15" = "Fait"
    "A16" = "Inspecteur";                                       "D16" = "A faire"
                      "B17" = "Déplacement";                    "D17" = "A faire"
                      "B18" = "Capture";                        "D18" = "A faire"
    "A19" = "Jeu";                                               "D19" = "A faire"
                      "B20" = "Start";                          "D20" = "A faire"
                      "B21" = "End";                            "D21" = "A faire"
    "A22" = "Jeton";                                             "D22" = "Fait"
}

# Row blocks (merged-row groupings) processed top to bottom; within each
# block the columns are filled in the order A, B, C, D.
$blocks = @(
    @(1),
    @(2),
    @(3),
    @(4),
    @(5, 6),
    @(7, 8),
    @(9, 10, 11),
    @(12, 13),
    @(14, 15),
    @(16, 17, 18),
    @(19, 20, 21),
    @(22)
)

$columnInfo = @(
    @{ Letter = "A"; Num = 1 },
    @{ Letter = "B"; Num = 2 },
    @{ Letter = "C"; Num = 3 },
    @{ Letter = "D"; Num = 4 }
)

foreach ($block in $blocks) {
    foreach ($colInfo in $columnInfo) {
        $col = $colInfo.Letter
        $colNum = $colInfo.Num
        foreach ($rowNum in $block) {
            $addr = "$col$rowNum"
            $cell = $ws.Cells.Item($rowNum, $colNum)
            if ($values.ContainsKey($addr)) {
                $cell.Value = $values[$addr]
            }
            if ($col -eq "A" -and $rowNum -ge 5 -and $rowNum -le 21) {
                $cell.HorizontalAlignment = -4108
                $cell.VerticalAlignment = -4108
            } else {
                $cell.HorizontalAlignment = -4108
            }
        }
    }
}

# ----------------------------------------------------------------------
# Merge the "Classe" column cells that share the same class grouping.
# Done after the values/styles are applied (only the top-left cell of
# each group carries a value, so nothing is lost by merging now).
# ----------------------------------------------------------------------
$null = $ws.Range("A5:A6").Merge()
$null = $ws.Range("A7:A8").Merge()
$null = $ws.Range("A9:A11").Merge()
$null = $ws.Range("A12:A13").Merge()
$null = $ws.Range("A14:A15").Merge()
$null = $ws.Range("A16:A18").Merge()
$null = $ws.Range("A19:A21").Merge()

# ----------------------------------------------------------------------
# Selection / view state
# ----------------------------------------------------------------------
$null = $ws.Range("G20").Select()
